# RPA datasets push 2024-07-25
# Updates the "02_38커뮤니케이션(최근일자기준)" sheet: a new IPO entry for
# "아이스크림미디어(구.시공미디어)" is inserted as row 4 (with updated
# bookbuilding dates/price range), and the subsequent rows (이엔셀, 엠83,
# 티디에스팜) shift down by one, which also removes the old duplicate
# "아이스크림미디어(구.시공미디어)" row that used to sit lower in the table.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(2)

$ws.Range("A4").Value = "아이스크림미디어(구.시공미디어)"
$ws.Range("B4").Value = "2024.08.09~08.16"
$ws.Range("C4").Value = "32,000~40,200"
$ws.Range("D4").Value = "-"
$ws.Range("E4").Value = 78720
$ws.Range("F4").Value = "삼성증권"

$ws.Range("A5").Value = "이엔셀"
$ws.Range("B5").Value = "2024.08.02~08.08"
$ws.Range("C5").Value = "13,600~15,300"
$ws.Range("D5").Value = "-"
$ws.Range("E5").Value = 21308
$ws.Range("F5").Value = "NH투자증권"

$ws.Range("A6").Value = "엠83"
$ws.Range("B6").Value = "2024.08.01~08.07"
$ws.Range("C6").Value = "11,000~13,000"
$ws.Range("D6").Value = "-"
$ws.Range("E6").Value = 16500
$ws.Range("F6").Value = "신영증권,유진투자증권"

$ws.Range("A7").Value = "티디에스팜"
$ws.Range("B7").Value = "2024.07.31~08.06"
$ws.Range("C7").Value = "9,500~10,700"
$ws.Range("D7").Value = "-"
$ws.Range("E7").Value = 9500
$ws.Range("F7").Value = "한국투자증권"
